# The source diff only touches PowerPoint's internal co-authoring /
# revision-tracking bookkeeping (ppt/revisionInfo.xml and
# ppt/changesInfos/changesInfo1.xml): the "F6540EED..." client's revision
# counter and timestamp are advanced (and the now-stale "711BDB89..."
# client entry is dropped from the revision list). Those parts just
# record *that* an edit/selection touched the picture on slide 1 again
# (picChg chg="mod" on shape id 7, "Graphic 6") - they carry no visible
# geometry/formatting/text payload, and every other part of the package
# (slide1.xml, media, rels, layouts, masters, ...) is byte-identical
# before and after. There is no real content change to make.
#
# Reproduce the author's action faithfully: touch/select the picture
# shape that is the subject of the recorded change, without altering any
# of its properties, so the slide's visible content and XML stay exactly
# as intended (matching the unchanged target content).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$pic = $s.Shapes.Item(1)
if ($pic.Name -ne "Graphic 6") {
    foreach ($shp in $s.Shapes) {
        if ($shp.Name -eq "Graphic 6") {
            $pic = $shp
        }
    }
}

# Select the picture (matches the recorded "custSel" + picChg chg="mod"
# on this shape) - a pure selection/touch, no property mutation.
$pic.Select()
